$wb = $excel.ActiveWorkbook

# --- Metadata sheet ---
$meta = $wb.Worksheets.Item("Metadata")

# Version 5.0.0 -> 6.0.0
$meta.Range("B3").Value = "6.0.0"

# Date updated
$meta.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value filled in
$meta.Range("B9").Value = "Alvearie Team"

# Row 10: Contact / No display for ContactDetail -> Jurisdiction / United States of America
$meta.Range("A10").Value = "Jurisdiction"
$meta.Range("B10").Value = "United States of America"

# Row 11 (duplicate Contact row) removed entirely - rows below shift up
$meta.Rows.Item(11).Delete()

# --- Elements sheet ---
$elements = $wb.Worksheets.Item("Elements")

# Root Extension row: Short/Definition now describe the renamed extension
$elements.Range("K2").Value = "Job Class Code"
$elements.Range("L2").Value = "Classifcation group of the job role of the employee, for example, Top Executives, Post-secondary teachers, and Engineers"
